# Updated cryptos list on Mon Apr 10 19:24:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.299.01"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.901.59"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.86%  "

# Row 5 - BNB
$ws.Range("D5").Value = "315.39"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.94%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3928"
$ws.Range("E8").Value = "  -1.28%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.08456"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10 - OKB
$ws.Range("D10").Value = "42.55"
$ws.Range("E10").Value = "  +1.49%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  +0.23%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "6.265"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.902.32"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14 - Solana
$ws.Range("D14").Value = "20.70"
$ws.Range("E14").Value = "  +0.23%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.338"
$ws.Range("E15").Value = "  +0.62%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.98%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "93.36"
$ws.Range("E17").Value = "  +1.94%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06738"
$ws.Range("E19").Value = "  -0.41%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +0.42%  "

# Row 21 - Dai
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -1.10%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.032"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "29.300.93"
$ws.Range("E23").Value = "  +2.32%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -0.28%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.220"
$ws.Range("E25").Value = "  -2.58%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.119.37"
$ws.Range("E26").Value = "  +1.02%  "

# Row 27 - was Monero, now EthereumClassic
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.02"
$ws.Range("E27").Value = "  +0.54%  "

# Row 28 - was EthereumClassic, now Monero
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "159.48"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.450"
$ws.Range("E29").Value = "  +2.20%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "128.26"
$ws.Range("E30").Value = "  +0.74%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "1.061"
$ws.Range("E31").Value = "  +0.73%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.1048"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "6.174"
$ws.Range("E33").Value = "  +6.09%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "3.668"
$ws.Range("E34").Value = "  +1.28%  "

# Row 35 - VeChain
$ws.Range("D35").Value = "0.02483"
$ws.Range("E35").Value = "  +1.58%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "0.06568"
$ws.Range("E36").Value = "  +0.87%  "

# Row 37 - FraxShare
$ws.Range("D37").Value = "9.083"
$ws.Range("E37").Value = "  +1.37%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "0.2199"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +2.87%  "

# Row 40 - InternetComputer(DFINITY)
$ws.Range("D40").Value = "5.153"
$ws.Range("E40").Value = "  +1.81%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "0.6503"
$ws.Range("E41").Value = "  +0.61%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.235"
$ws.Range("E42").Value = "  -2.48%  "

# Row 43 - Aptos
$ws.Range("E43").Value = "  +0.11%  "

# Row 44 - Decentraland
$ws.Range("D44").Value = "0.6069"
$ws.Range("E44").Value = "  -0.35%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "13.16"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46 - PancakeSwap
$ws.Range("D46").Value = "3.682"
$ws.Range("E46").Value = "  -1.07%  "

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "2.055"
$ws.Range("E47").Value = "  +2.30%  "

# Row 48 - EOS
$ws.Range("D48").Value = "1.228"
$ws.Range("E48").Value = "  +1.22%  "

# Row 49 - Quant
$ws.Range("D49").Value = "123.24"
$ws.Range("E49").Value = "  +0.26%  "

# Row 50 - WEMIXTOKEN
$ws.Range("E50").Value = "  -2.41%  "

# Row 51 - Aave
$ws.Range("D51").Value = "77.77"
$ws.Range("E51").Value = "  +0.65%  "
